$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" banner in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 23:42"

# Country data rows (row 4 through 204): country name + 7 numeric columns (B..H)
$data = @(
    ,@("Estados Unidos", 83144, 14933, 1864, 80079, 2112, 174, 1201)
    ,@("China", 81285, 0, 74051, 3947, 1235, 0, 3287)
    ,@("Italia", 80589, 6203, 10361, 62013, 3612, 712, 8215)
    ,@("España", 57786, 8271, 7015, 46406, 3166, 718, 4365)
    ,@("Alemania", 43938, 6615, 5673, 37998, 23, 61, 267)
    ,@("Iran", 29406, 2389, 10457, 16715, 2746, 157, 2234)
    ,@("Francia", 29155, 3922, 4948, 22511, 3375, 365, 1696)
    ,@("Suiza", 11811, 914, 131, 11489, 141, 38, 191)
    ,@("Reino Unido", 11658, 2129, 135, 10945, 163, 115, 578)
    ,@("Corea del Sur", 9241, 104, 4144, 4966, 59, 5, 131)
    ,@("Paises Bajos", 7431, 1019, 3, 6994, 761, 78, 434)
    ,@("Austria", 6898, 1310, 112, 6737, 28, 18, 49)
    ,@("Belgica", 6235, 1298, 675, 5340, 605, 42, 220)
    ,@("Canada", 4043, 634, 226, 3778, 120, 3, 39)
    ,@("Turquia", 3629, 1196, 26, 3528, 136, 16, 75)
    ,@("Portugal", 3544, 549, 43, 3441, 61, 17, 60)
    ,@("Noruega", 3369, 285, 6, 3349, 70, 0, 14)
    ,@("Australia", 2991, 315, 170, 2808, 11, 2, 13)
    ,@("Brasil", 2915, 361, 6, 2832, 18, 18, 77)
    ,@("Suecia", 2840, 314, 16, 2747, 176, 15, 77)
    ,@("Israel", 2693, 324, 68, 2617, 46, 3, 8)
    ,@("Malasia", 2031, 235, 215, 1793, 45, 3, 23)
    ,@("Chequia", 1925, 271, 10, 1906, 34, 3, 9)
    ,@("Dinamarca", 1877, 153, 1, 1835, 94, 7, 41)
    ,@("Irlanda", 1819, 255, 5, 1795, 47, 10, 19)
    ,@("Luxemburgo", 1453, 120, 6, 1438, 3, 1, 9)
    ,@("Ecuador", 1403, 192, 3, 1366, 58, 5, 34)
    ,@("Japon", 1399, 92, 359, 993, 57, 2, 47)
    ,@("Chile", 1306, 164, 22, 1280, 7, 1, 4)
    ,@("Polonia", 1221, 170, 7, 1198, 3, 2, 16)
    ,@("Pakistan", 1179, 116, 21, 1149, 5, 1, 9)
    ,@("Tailandia", 1045, 111, 88, 953, 4, 0, 4)
    ,@("Rumania", 1029, 123, 94, 912, 29, 6, 23)
    ,@("Arabia Saudita", 1012, 112, 33, 976, 6, 1, 3)
    ,@("Finlandia", 958, 78, 10, 943, 24, 2, 5)
    ,@("Sudafrica", 927, 218, 12, 915, 2, 0, 0)
    ,@("Indonesia", 893, 103, 35, 780, 0, 20, 78)
    ,@("Grecia", 892, 71, 36, 830, 53, 4, 26)
    ,@("Rusia", 840, 182, 38, 799, 8, 0, 3)
    ,@("Islandia", 802, 65, 68, 732, 11, 0, 2)
    ,@("India", 722, 65, 45, 661, 0, 4, 16)
    ,@("Crucero", 712, 0, 597, 105, 15, 0, 10)
    ,@("Filipinas", 707, 71, 28, 634, 1, 7, 45)
    ,@("Singapur", 683, 52, 172, 509, 17, 0, 2)
    ,@("Peru", 580, 100, 14, 557, 14, 0, 9)
    ,@("Eslovenia", 562, 34, 10, 546, 14, 1, 6)
    ,@("Panama", 558, 0, 2, 548, 20, 0, 8)
    ,@("Catar", 549, 12, 43, 506, 6, 0, 0)
    ,@("Estonia", 538, 134, 8, 529, 6, 0, 1)
    ,@("Argentina", 502, 0, 63, 430, 0, 1, 9)
    ,@("Croacia", 495, 53, 22, 470, 14, 2, 3)
    ,@("Egipto", 495, 39, 102, 369, 0, 3, 24)
    ,@("Colombia", 491, 21, 8, 477, 0, 2, 6)
    ,@("Republica Dominicana", 488, 96, 3, 475, 0, 0, 10)
    ,@("Mexico", 475, 70, 4, 465, 1, 1, 6)
    ,@("Barein", 458, 39, 204, 250, 1, 0, 4)
    ,@("Serbia", 457, 73, 15, 435, 21, 3, 7)
    ,@("Hong Kong", 453, 42, 110, 339, 5, 0, 4)
    ,@("Irak", 382, 36, 105, 241, 0, 7, 36)
    ,@("Libano", 368, 35, 23, 339, 3, 0, 6)
    ,@("Argelia", 367, 65, 29, 313, 0, 4, 25)
    ,@("Emiratos Arabes Unidos", 333, 0, 52, 279, 2, 0, 2)
    ,@("Lituania", 299, 25, 1, 294, 1, 0, 4)
    ,@("Armenia", 290, 25, 18, 271, 6, 1, 1)
    ,@("Nueva Zelanda", 283, 0, 27, 256, 0, 0, 0)
    ,@("Marruecos", 275, 50, 8, 256, 1, 5, 11)
    ,@("Bulgaria", 264, 22, 8, 253, 8, 0, 3)
    ,@("Hungria", 261, 35, 28, 223, 6, 0, 10)
    ,@("Taiwan", 252, 17, 29, 221, 0, 0, 2)
    ,@("Letonia", 244, 23, 1, 243, 0, 0, 0)
    ,@("Costa Rica", 231, 30, 2, 227, 5, 0, 2)
    ,@("Eslovaquia", 226, 10, 2, 224, 2, 0, 0)
    ,@("Principado de Andorra", 224, 36, 1, 220, 6, 2, 3)
    ,@("Uruguay", 217, 0, 0, 217, 3, 0, 0)
    ,@("Jordania", 212, 40, 1, 211, 0, 0, 0)
    ,@("San Marino", 208, 0, 4, 183, 12, 0, 21)
    ,@("Kuwait", 208, 13, 49, 159, 7, 0, 0)
    ,@("Republica de Macedonia", 201, 24, 3, 195, 1, 0, 3)
    ,@("Tunez", 197, 24, 2, 189, 10, 1, 6)
    ,@("Ucrania", 196, 51, 1, 190, 0, 0, 5)
    ,@("Bosnia y Herzegovina", 191, 15, 2, 186, 1, 0, 3)
    ,@("Moldavia", 177, 28, 2, 174, 28, 0, 1)
    ,@("Albania", 174, 28, 17, 151, 3, 1, 6)
    ,@("Vietnam", 153, 5, 20, 133, 3, 0, 0)
    ,@("Burkina Faso", 152, 6, 10, 135, 0, 3, 7)
    ,@("Republica de Chipre", 146, 14, 4, 139, 3, 0, 3)
    ,@("Islas Feroe", 140, 8, 47, 93, 0, 0, 0)
    ,@("Reunion", 135, 24, 1, 134, 0, 0, 0)
    ,@("Malta", 134, 5, 2, 132, 1, 0, 0)
    ,@("Ghana", 132, 64, 1, 127, 0, 0, 4)
    ,@("Azerbaiyan", 122, 29, 15, 104, 6, 1, 3)
    ,@("Brunei", 114, 5, 5, 109, 1, 0, 0)
    ,@("Kazajistan", 111, 30, 2, 108, 0, 1, 1)
    ,@("Oman", 109, 10, 23, 86, 0, 0, 0)
    ,@("Venezuela", 107, 1, 15, 92, 2, 0, 0)
    ,@("Sri Lanka", 106, 4, 7, 99, 3, 0, 0)
    ,@("Senegal", 105, 6, 9, 96, 0, 0, 0)
    ,@("Costa de Marfil", 96, 16, 3, 93, 0, 0, 0)
    ,@("Camboya", 96, 0, 10, 86, 1, 0, 0)
    ,@("Afganistan", 94, 10, 2, 88, 0, 2, 4)
    ,@("Bielorrusia", 86, 0, 29, 57, 2, 0, 0)
    ,@("Estado de Palestina", 84, 13, 17, 66, 0, 0, 1)
    ,@("Mauricio", 81, 33, 0, 79, 1, 0, 2)
    ,@("Georgia", 79, 4, 11, 68, 1, 0, 0)
    ,@("Uzbekistan", 75, 15, 0, 75, 4, 0, 0)
    ,@("Camerun", 75, 0, 2, 72, 0, 0, 1)
    ,@("Guadalupe", 73, 0, 0, 72, 4, 0, 1)
    ,@("Montenegro", 69, 16, 0, 68, 1, 0, 1)
    ,@("Cuba", 67, 10, 1, 65, 2, 0, 1)
    ,@("Martinica", 66, 0, 0, 65, 7, 0, 1)
    ,@("Trinidad yTobago", 65, 5, 0, 64, 0, 0, 1)
    ,@("Nigeria", 65, 14, 2, 62, 0, 0, 1)
    ,@("Liechtenstein", 56, 5, 0, 56, 0, 0, 0)
    ,@("Honduras", 52, 0, 0, 51, 0, 1, 1)
    ,@("Consejo Danes para los Refugiados", 51, 3, 0, 48, 0, 1, 3)
    ,@("Kirguistan", 44, 0, 0, 44, 0, 0, 0)
    ,@("Banglades", 44, 5, 11, 28, 1, 0, 5)
    ,@("Bolivia", 43, 11, 0, 43, 0, 0, 0)
    ,@("Ruanda", 41, 0, 0, 41, 0, 0, 0)
    ,@("Paraguay", 41, 4, 0, 38, 1, 0, 3)
    ,@("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
    ,@("Mayotte", 36, 0, 0, 36, 0, 0, 0)
    ,@("Gibraltar", 35, 9, 13, 22, 0, 0, 0)
    ,@("Monaco", 33, 2, 1, 32, 0, 0, 0)
    ,@("Guam", 32, 0, 0, 31, 0, 0, 1)
    ,@("Kenia", 31, 3, 1, 29, 0, 1, 1)
    ,@("Macao", 31, 0, 10, 21, 0, 0, 0)
    ,@("Polinesia Francesa", 30, 5, 0, 30, 0, 0, 0)
    ,@("Guayana Francesa", 28, 0, 6, 22, 0, 0, 0)
    ,@("Jamaica", 26, 0, 2, 23, 0, 0, 1)
    ,@("Isla de Man", 25, 2, 0, 25, 0, 0, 0)
    ,@("Guatemala", 24, 0, 4, 19, 0, 0, 1)
    ,@("Madagascar", 23, 4, 0, 23, 0, 0, 0)
    ,@("Togo", 23, 0, 1, 22, 0, 0, 0)
    ,@("Aruba", 19, 0, 1, 18, 0, 0, 0)
    ,@("Barbados", 18, 0, 0, 18, 0, 0, 0)
    ,@("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
    ,@("Zambia", 16, 4, 0, 16, 0, 0, 0)
    ,@("Bermudas", 15, 8, 2, 13, 0, 0, 0)
    ,@("Uganda", 14, 0, 0, 14, 0, 0, 0)
    ,@("Nueva Caledonia", 14, 0, 0, 14, 0, 0, 0)
    ,@("Tanzania", 13, 0, 0, 13, 0, 0, 0)
    ,@("El Salvador", 13, 4, 0, 13, 0, 0, 0)
    ,@("Maldivas", 13, 0, 8, 5, 0, 0, 0)
    ,@("Guinea Ecuatorial", 12, 3, 0, 12, 0, 0, 0)
    ,@("Etiopia", 12, 0, 0, 12, 0, 0, 0)
    ,@("Dominica", 11, 0, 0, 11, 0, 0, 0)
    ,@("Republica de Yibuti", 11, 0, 0, 11, 0, 0, 0)
    ,@("Mongolia", 11, 1, 0, 11, 0, 0, 0)
    ,@("San Martin (Parte Francesa)", 11, 0, 0, 11, 0, 0, 0)
    ,@("Niger", 10, 3, 0, 9, 0, 0, 1)
    ,@("Bahamas", 9, 4, 1, 8, 0, 0, 0)
    ,@("Haiti", 8, 0, 0, 8, 0, 0, 0)
    ,@("Surinam", 8, 0, 0, 8, 0, 0, 0)
    ,@("Islas Caimanes", 8, 0, 0, 7, 0, 0, 1)
    ,@("Namibia", 8, 1, 2, 6, 0, 0, 0)
    ,@("Mozambique", 7, 2, 0, 7, 0, 0, 0)
    ,@("Antigua y Barbuda", 7, 4, 0, 7, 0, 0, 0)
    ,@("Seychelles", 7, 0, 0, 7, 0, 0, 0)
    ,@("Granada", 7, 6, 0, 7, 0, 0, 0)
    ,@("Gabon", 7, 1, 0, 6, 0, 0, 1)
    ,@("Laos", 6, 3, 0, 6, 0, 0, 0)
    ,@("Suazilandia", 6, 2, 0, 6, 0, 0, 0)
    ,@("Eritrea", 6, 2, 0, 6, 0, 0, 0)
    ,@("Benin", 6, 0, 0, 6, 0, 0, 0)
    ,@("Groenlandia", 6, 0, 2, 4, 0, 0, 0)
    ,@("Curazao", 6, 0, 2, 3, 0, 0, 1)
    ,@("Siria", 5, 0, 0, 5, 0, 0, 0)
    ,@("Fiyi", 5, 0, 0, 5, 0, 0, 0)
    ,@("Guyana", 5, 0, 0, 4, 0, 0, 1)
    ,@("Congo", 4, 0, 0, 4, 0, 0, 0)
    ,@("Mali", 4, 2, 0, 4, 0, 0, 0)
    ,@("Santa Sede", 4, 0, 0, 4, 0, 0, 0)
    ,@("Angola", 4, 1, 0, 4, 0, 0, 0)
    ,@("Guinea", 4, 0, 0, 4, 0, 0, 0)
    ,@("Cabo Verde", 4, 0, 0, 3, 0, 0, 1)
    ,@("San Martin (Parte Holandesa)", 3, 0, 0, 3, 0, 0, 0)
    ,@("Mauritania", 3, 1, 0, 3, 0, 0, 0)
    ,@("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0)
    ,@("Liberia", 3, 0, 0, 3, 0, 0, 0)
    ,@("San Bartolome", 3, 0, 0, 3, 0, 0, 0)
    ,@("Birmania", 3, 0, 0, 3, 0, 0, 0)
    ,@("Republica del Chad", 3, 0, 0, 3, 0, 0, 0)
    ,@("Santa Lucia", 3, 0, 0, 3, 0, 0, 0)
    ,@("Gambia", 3, 0, 0, 2, 0, 0, 1)
    ,@("Sudan", 3, 0, 0, 2, 0, 0, 1)
    ,@("Zimbabue", 3, 0, 0, 2, 0, 0, 1)
    ,@("Nepal", 3, 0, 1, 2, 0, 0, 0)
    ,@("Guinea-Bisau", 2, 0, 0, 2, 0, 0, 0)
    ,@("San Cristobal y Nieves", 2, 0, 0, 2, 0, 0, 0)
    ,@("Islas Turcas y Caicos", 2, 1, 0, 2, 0, 0, 0)
    ,@("Belice", 2, 0, 0, 2, 0, 0, 0)
    ,@("Nicaragua", 2, 0, 0, 2, 0, 0, 0)
    ,@("Butan", 2, 0, 0, 2, 0, 0, 0)
    ,@("Islas Virgenes Britanicas", 2, 0, 0, 2, 0, 0, 0)
    ,@("Montserrat", 2, 1, 0, 2, 0, 0, 0)
    ,@("Somalia", 2, 1, 0, 2, 0, 0, 0)
    ,@("San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0)
    ,@("Timor Oriental", 1, 0, 0, 1, 0, 0, 0)
    ,@("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0)
    ,@("Libia", 1, 0, 0, 1, 0, 0, 0)
)

$row = 4
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = [double]$item[1]
    $ws.Cells.Item($row, 3).Value = [double]$item[2]
    $ws.Cells.Item($row, 4).Value = [double]$item[3]
    $ws.Cells.Item($row, 5).Value = [double]$item[4]
    $ws.Cells.Item($row, 6).Value = [double]$item[5]
    $ws.Cells.Item($row, 7).Value = [double]$item[6]
    $ws.Cells.Item($row, 8).Value = [double]$item[7]
    $row = $row + 1
}

Write-Host "Updated $($row - 4) country rows"
